$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 894.8570999999999
$ws.Range("I32").Value = 826.3333
$ws.Range("J32").Value = 946.25
$ws.Range("K32").Value = 826.3333
$ws.Range("L32").Value = 946.25
$ws.Range("M32").Value = -500.3333
$ws.Range("N32").Value = -1598.25
$ws.Range("H70").Value = 1306.9714
$ws.Range("J70").Value = 1277.7646
$ws.Range("L70").Value = 3833.2938
$ws.Range("N70").Value = -4373.293799999999
$ws.Range("H73").Value = 1306.9714
$ws.Range("J73").Value = 1277.7646
$ws.Range("L73").Value = 3833.2938
$ws.Range("N73").Value = -5705.293799999999
$ws.Range("H86").Value = 5733.0835
$ws.Range("I86").Value = 5423.25
$ws.Range("J86").Value = 5888
$ws.Range("K86").Value = 5423.25
$ws.Range("L86").Value = 5888
$ws.Range("M86").Value = -4300.25
$ws.Range("N86").Value = -8134
$ws.Range("H88").Value = 24189.25
$ws.Range("I88").Value = 71957.125
$ws.Range("J88").Value = 5082.1
$ws.Range("K88").Value = 71957.125
$ws.Range("L88").Value = 5082.1
$ws.Range("M88").Value = -71551.125
$ws.Range("N88").Value = -5894.1
$ws.Range("H89").Value = 5733.0835
$ws.Range("I89").Value = 5423.25
$ws.Range("J89").Value = 5888
$ws.Range("K89").Value = 27116.25
$ws.Range("L89").Value = 29440
$ws.Range("M89").Value = -21500.25
$ws.Range("N89").Value = -40672
$ws.Range("H91").Value = 24189.25
$ws.Range("I91").Value = 71957.125
$ws.Range("J91").Value = 5082.1
$ws.Range("K91").Value = 71957.125
$ws.Range("L91").Value = 5082.1
$ws.Range("M91").Value = -70553.125
$ws.Range("N91").Value = -7890.1
$ws.Range("H113").Value = 4620.2
$ws.Range("I113").Value = 4080.9092
$ws.Range("J113").Value = 4932.421
$ws.Range("K113").Value = 4080.9092
$ws.Range("L113").Value = 4932.421
$ws.Range("M113").Value = -826.9092000000001
$ws.Range("N113").Value = -11440.421
$ws.Range("H116").Value = 3100.8
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3100.8
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 3100.8
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -9984.799999999999
$ws.Range("H118").Value = 532.7059
$ws.Range("I118").Value = 370.4
$ws.Range("J118").Value = 1750
$ws.Range("K118").Value = 1111.2
$ws.Range("L118").Value = 5250
$ws.Range("M118").Value = 545.8000000000002
$ws.Range("N118").Value = -8564
$ws.Range("H137").Value = 1355.6809
$ws.Range("I137").Value = 1144.9
$ws.Range("J137").Value = 1727.6471
$ws.Range("K137").Value = 3434.7
$ws.Range("L137").Value = 5182.9413
$ws.Range("M137").Value = -884.7000000000003
$ws.Range("N137").Value = -10282.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31748
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -98736
$ws.Range("H88").Value = 1933
$ws.Range("J88").Value = 1799.3334
$ws.Range("L88").Value = 1799.3334
$ws.Range("N88").Value = -2611.3334
$ws.Range("H91").Value = 1933
$ws.Range("J91").Value = 1799.3334
$ws.Range("L91").Value = 1799.3334
$ws.Range("N91").Value = -4607.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2043.3043
$ws.Range("I86").Value = 1965.5555
$ws.Range("J86").Value = 2093.2856
$ws.Range("K86").Value = 1965.5555
$ws.Range("L86").Value = 2093.2856
$ws.Range("M86").Value = -842.5554999999999
$ws.Range("N86").Value = -4339.2856
$ws.Range("H89").Value = 2043.3043
$ws.Range("I89").Value = 1965.5555
$ws.Range("J89").Value = 2093.2856
$ws.Range("K89").Value = 9827.7775
$ws.Range("L89").Value = 10466.428
$ws.Range("M89").Value = -4211.7775
$ws.Range("N89").Value = -21698.428
$ws.Range("H117").Value = 31000
$ws.Range("J117").Value = 31000
$ws.Range("L117").Value = 31000
$ws.Range("N117").Value = -40178

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 6.375
$ws.Range("I14").Value = 6.375
$ws.Range("K14").Value = 19.125
$ws.Range("M14").Value = 153.875
$ws.Range("H23").Value = 173.4762
$ws.Range("I23").Value = 59.88889
$ws.Range("J23").Value = 258.66666
$ws.Range("K23").Value = 179.66667
$ws.Range("L23").Value = 775.9999799999999
$ws.Range("M23").Value = 55.33332999999999
$ws.Range("N23").Value = -1245.99998
$ws.Range("H68").Value = 1533.3334
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 1533.3334
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H86").Value = 450
$ws.Range("J86").Value = 450
$ws.Range("L86").Value = 1350
$ws.Range("N86").Value = -3722
$ws.Range("H89").Value = 450
$ws.Range("J89").Value = 450
$ws.Range("L89").Value = 4050
$ws.Range("N89").Value = -15906

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 50794.05
$ws.Range("I22").Value = 200476.2
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 200476.2
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -200181.2
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 50794.05
$ws.Range("I27").Value = 200476.2
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 200476.2
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -200369.2
$ws.Range("N27").Value = -1114
$ws.Range("H82").Value = 1450.4445
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1450.4445
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1450.4445
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2172.4445
$ws.Range("H85").Value = 1450.4445
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1450.4445
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1450.4445
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3946.4445
$ws.Range("H122").Value = 4510.9375
$ws.Range("I122").Value = 5612.5
$ws.Range("J122").Value = 3850
$ws.Range("K122").Value = 16837.5
$ws.Range("L122").Value = 11550
$ws.Range("M122").Value = -14387.5
$ws.Range("N122").Value = -16450

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6565.2
$ws.Range("J74").Value = 6565.2
$ws.Range("L74").Value = 6565.2
$ws.Range("N74").Value = -8437.200000000001
$ws.Range("H77").Value = 6565.2
$ws.Range("J77").Value = 6565.2
$ws.Range("L77").Value = 19695.6
$ws.Range("N77").Value = -29055.6
$ws.Range("H81").Value = 1533.8
$ws.Range("I81").Value = 1496.3334
$ws.Range("J81").Value = 1590
$ws.Range("K81").Value = 2992.6668
$ws.Range("L81").Value = 3180
$ws.Range("M81").Value = -1931.6668
$ws.Range("N81").Value = -5302
$ws.Range("H84").Value = 1533.8
$ws.Range("I84").Value = 1496.3334
$ws.Range("J84").Value = 1590
$ws.Range("K84").Value = 14963.334
$ws.Range("M84").Value = -9659.333999999999
$ws.Range("N84").Value = -26508
